# Generate Report for Handback
# Updates the localization-status report to reflect that the
# 9bfdd7b6-0537-41ef-9a45-339b14cae5cd.md file has now been handed back
# (instead of merely "Ready for handoff") for both the zh-cn and de-de
# target languages, and records the handback timestamps.

$wb = $excel.ActiveWorkbook

$handedBack = "Handed back: in sync with en-US"

# --- Overview sheet: update the per-language status columns for row 3 ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $handedBack
$overview.Range("C3").Value = $handedBack

# --- zh-cn sheet: update Status and Latest Handback DateTime for row 3 ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $handedBack
$zhcn.Range("H3").Value = "2016-03-18 04:23:28"

# --- de-de sheet: update Status and Latest Handback DateTime for row 3 ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $handedBack
$dede.Range("H3").Value = "2016-03-18 04:23:32"
